# removed cleaning on table sP1
# Swaps the two table titles ("Sector 29" <-> "Sector 13") and refreshes all
# summary-statistic values for both tables on the active sheet (re-run of the
# analysis without the data-cleaning step).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table headers (title swap) -------------------------------------------
$ws.Range("A1").Value = "Sector 13"
$ws.Range("A8").Value = "Sector 29"

# --- Table 1 (rows 2-6) -----------------------------------------------------
$table1 = @(
    @(132.92127965492452, 147.71975601491019, 1008.9617118503314, 1063.504916237069, -14.798476359985671, 27.413990712196949),
    @(71848.445929923197, 58150.946746828689, 1087504.9305938559, 759169.99574361613, 13697.499183094507, 24673.161334857159),
    @(3178.5482154052174, 3074.7134027941652, 28886.462084369377, 25037.785176689635, 103.83481261105226, 708.84694584384204),
    @(46754.904746624074, 46125.585812180259, 724089.21607989911, 647590.91457900731, 629.31893444381421, 18308.404903139359),
    @(19456.589893134285, 18257.890227933472, 249812.66980913674, 227306.12426804498, 1198.6996652008129, 6443.892666909439)
)

for ($i = 0; $i -lt $table1.Length; $i++) {
    $row = 2 + $i
    $values = $table1[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}

# --- Table 2 (rows 9-13) -----------------------------------------------------
$table2 = @(
    @(20.112409479921002, 23.872162485065711, 41.511698434113832, 49.059182334552915, -3.7597530051447094, 0.80837371229434918),
    @(3708.0696605941039, 3975.713398269937, 11230.372591117222, 10107.992072367198, -267.64373767583311, 189.1118745961534),
    @(344.35838964046087, 429.38602486757401, 1593.1220013942445, 1612.6354098410256, -85.027635227113137, 28.113523705476979),
    @(1648.972774003666, 2055.0905018921512, 6164.5965267078182, 6271.8004843207045, -406.11772788848521, 111.7965765741109),
    @(2017.4392632041686, 2104.7503877541872, 5258.8515317592301, 4983.3053038042881, -87.311124550018576, 93.291310057966655)
)

for ($i = 0; $i -lt $table2.Length; $i++) {
    $row = 9 + $i
    $values = $table2[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
